# Update the "Datos actualizados" timestamp caption in A1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Datos actualizados a 9 de Abril de 2020 a las 12:52"

# Update the rows whose case-count figures changed (and, where the update
# caused a province to overtake its neighbour in the ranking, swap the two
# rows so the table stays sorted by "Casos totales" descending).

$ws.Range("A4").Value  = "Madrid"
$ws.Range("B4").Value  = 43877
$ws.Range("C4").Value  = 21121
$ws.Range("D4").Value  = 16956
$ws.Range("E4").Value  = 5800

$ws.Range("A10").Value = "Navarra"
$ws.Range("B10").Value = 3575
$ws.Range("C10").Value = 462
$ws.Range("D10").Value = 2899
$ws.Range("E10").Value = 214

$ws.Range("A15").Value = "Zaragoza"
$ws.Range("B15").Value = 2788
$ws.Range("C15").Value = 534
$ws.Range("D15").Value = 1967
$ws.Range("E15").Value = 287

$ws.Range("A16").Value = "Castilla-La Mancha"
$ws.Range("B16").Value = 2780
$ws.Range("C16").Value = 71
$ws.Range("D16").Value = 2446
$ws.Range("E16").Value = 263

$ws.Range("A28").Value = "Caceres"
$ws.Range("B28").Value = 1513
$ws.Range("C28").Value = 177
$ws.Range("D28").Value = 1101
$ws.Range("E28").Value = 235

$ws.Range("A29").Value = "Segovia"
$ws.Range("B29").Value = 1480
$ws.Range("C29").Value = 415
$ws.Range("D29").Value = 943
$ws.Range("E29").Value = 122

$ws.Range("A42").Value = "Badajoz"
$ws.Range("B42").Value = 760
$ws.Range("C42").Value = 204
$ws.Range("D42").Value = 508
$ws.Range("E42").Value = 48

$ws.Range("A43").Value = "Ourense"
$ws.Range("B43").Value = 751
$ws.Range("C43").Value = 333
$ws.Range("D43").Value = 660
$ws.Range("E43").Value = 22

$ws.Range("A47").Value = "Huesca"
$ws.Range("B47").Value = 444
$ws.Range("C47").Value = 74
$ws.Range("D47").Value = 312
$ws.Range("E47").Value = 58

$ws.Range("A48").Value = "Gran Canaria"
$ws.Range("B48").Value = 434
$ws.Range("C48").Value = 104
$ws.Range("D48").Value = 305
$ws.Range("E48").Value = 25

$ws.Range("A49").Value = "Teruel"
$ws.Range("B49").Value = 416
$ws.Range("C49").Value = 101
$ws.Range("D49").Value = 276
$ws.Range("E49").Value = 39
